$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Privat")

$data = @(
    @("Dr. Schmidt", "Wednesday", "2023-11-15", "14:00"),
    @("Dr. Schmidt", "Friday",    "2023-11-17", "14:30"),
    @("Dr. Schmidt", "Friday",    "2023-11-17", "14:30")
)

$startRow = 7
$endRow = $startRow + $data.Count - 1

# Columns C (Date) and D (Heure) hold date/time-looking text such as
# "2023-11-15" and "14:00". Temporarily force text format on just those
# columns so Excel doesn't auto-convert the strings into numeric
# date/time values, then restore the default "Normal" style so no extra
# formatting is left behind on the new cells (matches the source rows,
# which carry no explicit style).
$dateTimeRange = $ws.Range("C$startRow`:D$endRow")
$dateTimeRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

$dateTimeRange.Style = "Normal"
